$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original formatting on the Price/Volume columns, then force them to
# Text so Excel does not reinterpret numeric-looking strings (e.g. "1.00") as
# numbers and silently drop the formatting (trailing zeros, thousands dots, etc).
$priceVolRange = $ws.Range("D2:E51")
$origStyle = $priceVolRange.Style
$priceVolRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.576.01"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3
$ws.Range("D3").Value = "3.277.23"
$ws.Range("E3").Value = "  +3.88%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "596.93"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "142.06"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "3.267.94"
$ws.Range("E8").Value = "  +3.85%  "

# Row 9
$ws.Range("E9").Value = "  -0.48%  "

# Row 10
$ws.Range("E10").Value = "  -0.19%  "

# Row 11
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +0.40%  "

# Row 13
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -1.48%  "

# Row 14
$ws.Range("D14").Value = "35.01"
$ws.Range("E14").Value = "  -0.17%  "

# Row 15
$ws.Range("D15").Value = "3.810.59"
$ws.Range("E15").Value = "  +3.93%  "

# Row 16
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").Value = "3.259.85"
$ws.Range("E17").Value = "  +3.60%  "

# Row 18
$ws.Range("D18").Value = "63.600.92"
$ws.Range("E18").Value = "  -0.81%  "

# Row 19
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +0.01%  "

# Row 20
$ws.Range("D20").Value = "480.86"
$ws.Range("E20").Value = "  -1.71%  "

# Row 21
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  -2.62%  "

# Row 22
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +3.42%  "

# Row 23
$ws.Range("D23").Value = "8.01"
$ws.Range("E23").Value = "  +4.58%  "

# Row 24
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "13.52"
$ws.Range("E24").Value = "  +1.22%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "84.62"
$ws.Range("E25").Value = "  -3.21%  "

# Row 26
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +5.23%  "

# Row 29
$ws.Range("D29").Value = "8.20"
$ws.Range("E29").Value = "  -0.54%  "

# Row 30
$ws.Range("D30").Value = "2.17"
$ws.Range("E30").Value = "  +4.96%  "

# Row 31
$ws.Range("D31").Value = "28.20"
$ws.Range("E31").Value = "  +3.39%  "

# Row 32
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$ws.Range("E33").Value = "  -2.22%  "

# Row 34
$ws.Range("D34").Value = "2.60"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  -1.37%  "

# Row 36
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
$ws.Range("D37").Value = "53.07"
$ws.Range("E37").Value = "  +0.61%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0729"
$ws.Range("E38").Value = "  -2.40%  "

# Row 39
$ws.Range("D39").Value = "0.0397"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("D40").Value = "425.63"
$ws.Range("E40").Value = "  -3.00%  "

# Row 41
$ws.Range("D41").Value = "2.82"
$ws.Range("E41").Value = "  -5.74%  "

# Row 42
$ws.Range("D42").Value = "8.45"
$ws.Range("E42").Value = "  +1.90%  "

# Row 43
$ws.Range("D43").Value = "3.005.10"
$ws.Range("E43").Value = "  +3.10%  "

# Row 44
$ws.Range("D44").Value = "0.113"
$ws.Range("E44").Value = "  -6.23%  "

# Row 45
$ws.Range("D45").Value = "0.270"
$ws.Range("E45").Value = "  +3.95%  "

# Row 46
$ws.Range("E46").Value = "  -0.85%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "26.22"
$ws.Range("E47").Value = "  +0.73%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  -3.47%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.115"
$ws.Range("E50").Value = "  +0.49%  "

# Row 51
$ws.Range("D51").Value = "119.66"
$ws.Range("E51").Value = "  -0.64%  "

# Restore the original number formatting/style now that the text is committed.
$priceVolRange.Style = $origStyle
